$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.082.81"
$ws.Range("E2").Value = "  -4.38%  "
$ws.Range("D3").Value = "2.977.84"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.79%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +3.02%  "
$ws.Range("D9").Value = "2.970.37"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.86%  "
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "3.470.51"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("E17").Value = "  +6.48%  "
$ws.Range("D18").Value = "2.975.29"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").Value = "58.020.10"
$ws.Range("E19").Value = "  -4.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "421.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.689"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.54%  "
$ws.Range("E30").Value = "  +4.56%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.101"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.946"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").Value = "0.0₃0704"
$ws.Range("E37").Value = "  +4.30%  "
$ws.Range("E38").Value = "  -2.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "379.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").Value = "2.692.07"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.84%  "
$ws.Range("E48").Value = "  +2.46%  "
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("E51").Value = "  -0.49%  "
